$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date-serial that was bumped by one day
# (46061 -> 46062) for every data row on the sheet. Determine the data
# range from the sheet's used range (row 1 is the header row).
$used = $ws.UsedRange
$firstRow = 2
$lastRow = $used.Row + $used.Rows.Count - 1
if ($lastRow -lt $firstRow) { $lastRow = $firstRow }

$ws.Range("C" + $firstRow + ":C" + $lastRow).Value = 46062
